$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row/Col/Value triples for the updated vm_pu results (Case with 380 kV done)
$updates = @(
    @(2, 2, 1.02),
    @(2, 3, 1.03619328014384),
    @(2, 4, 1.038837636877475),
    @(2, 5, 1.043969571154188),
    @(2, 6, 1.052083561064594),
    @(2, 9, 1.034311452478081),
    @(2, 10, 1.041302932441949),
    @(2, 11, 1.041624532558426),
    @(2, 12, 1.046741939805645),
    @(2, 13, 1.054833274861561),
    @(2, 14, 1.042781702045133),
    @(3, 2, 1.02),
    @(3, 3, 1.037284424865084),
    @(3, 4, 1.039843006486156),
    @(3, 5, 1.044974626868057),
    @(3, 6, 1.053258413668583),
    @(3, 9, 1.034493273632572),
    @(3, 10, 1.042037225991341),
    @(3, 11, 1.042439482123166),
    @(3, 12, 1.047557631996732),
    @(3, 13, 1.055819972693901),
    @(3, 14, 1.043517038375589),
    @(4, 2, 1.02),
    @(4, 3, 1.037989996381771),
    @(4, 4, 1.0404933559343),
    @(4, 5, 1.045624985271439),
    @(4, 6, 1.054018923919444),
    @(4, 9, 1.03460859094797),
    @(4, 10, 1.042511341561719),
    @(4, 11, 1.042966010253636),
    @(4, 12, 1.04808485059762),
    @(4, 13, 1.05645815869345),
    @(4, 14, 1.043991827244459),
    @(5, 2, 1.02),
    @(5, 3, 1.038286506344037),
    @(5, 4, 1.040766717015726),
    @(5, 5, 1.045898400910227),
    @(5, 6, 1.054338715205929),
    @(5, 9, 1.034656511468768),
    @(5, 10, 1.042710415314785),
    @(5, 11, 1.043187171570005),
    @(5, 12, 1.048306352499645),
    @(5, 13, 1.056726386937613),
    @(5, 14, 1.044191183705088),
    @(6, 2, 1.02),
    @(6, 3, 1.038336285106453),
    @(6, 4, 1.040812612860083),
    @(6, 5, 1.045944308894222),
    @(6, 6, 1.054392413851098),
    @(6, 9, 1.034664524752554),
    @(6, 10, 1.042743826343109),
    @(6, 11, 1.043224294324352),
    @(6, 12, 1.048343535391664),
    @(6, 13, 1.05677141982879),
    @(6, 14, 1.044224642180905),
    @(7, 2, 1.02),
    @(7, 3, 1.037993958802805),
    @(7, 4, 1.040497008778395),
    @(7, 5, 1.04562863864492),
    @(7, 6, 1.054023196699742),
    @(7, 9, 1.034609233460756),
    @(7, 10, 1.042514002555302),
    @(7, 11, 1.042968966172794),
    @(7, 12, 1.048087810869941),
    @(7, 13, 1.056461743027995),
    @(7, 14, 1.043994492016958),
    @(8, 2, 1.02),
    @(8, 3, 1.036562135791543),
    @(8, 4, 1.039177446056582),
    @(8, 5, 1.044309230355009),
    @(8, 6, 1.052480545742606),
    @(8, 9, 1.034373382688175),
    @(8, 10, 1.041551301952056),
    @(8, 11, 1.04190011403132),
    @(8, 12, 1.047017728813251),
    @(8, 13, 1.055166790835235),
    @(8, 14, 1.043030424268432),
    @(9, 2, 1.02),
    @(9, 3, 1.034035432141881),
    @(9, 4, 1.036850731583285),
    @(9, 5, 1.041984411242319),
    @(9, 6, 1.049764489230544),
    @(9, 9, 1.03393992791603),
    @(9, 10, 1.039847071484809),
    @(9, 11, 1.04001053666039),
    @(9, 12, 1.045127589892894),
    @(9, 13, 1.052882809711441),
    @(9, 14, 1.041323773598444),
    @(10, 2, 1.02),
    @(10, 3, 1.032348460023266),
    @(10, 4, 1.035298582777509),
    @(10, 5, 1.040434616296027),
    @(10, 6, 1.047955291210362),
    @(10, 9, 1.033638968974135),
    @(10, 10, 1.038705642304311),
    @(10, 11, 1.038746689306265),
    @(10, 12, 1.043864444974722),
    @(10, 13, 1.051358711787146),
    @(10, 14, 1.040180723457586),
    @(11, 2, 1.02),
    @(11, 3, 1.031617376508227),
    @(11, 4, 1.034626242448789),
    @(11, 5, 1.039763553645304),
    @(11, 6, 1.047172236920136),
    @(11, 9, 1.033505809475482),
    @(11, 10, 1.03821013586956),
    @(11, 11, 1.038198446234347),
    @(11, 12, 1.043316760926306),
    @(11, 13, 1.05069841133536),
    @(11, 14, 1.039684513346862),
    @(12, 2, 1.02),
    @(12, 3, 1.031345726183878),
    @(12, 4, 1.034376467312337),
    @(12, 5, 1.039514291916965),
    @(12, 6, 1.046881425941313),
    @(12, 9, 1.033455921211334),
    @(12, 10, 1.038025892962242),
    @(12, 11, 1.037994655327529),
    @(12, 12, 1.043113215918548),
    @(12, 13, 1.050453092155521),
    @(12, 14, 1.039500008793483),
    @(13, 2, 1.02),
    @(13, 3, 1.03140400034436),
    @(13, 4, 1.034430046654759),
    @(13, 5, 1.039567759374928),
    @(13, 6, 1.046943803623934),
    @(13, 9, 1.033466641729468),
    @(13, 10, 1.038065422285362),
    @(13, 11, 1.038038375919771),
    @(13, 12, 1.043156882034173),
    @(13, 13, 1.05050571641615),
    @(13, 14, 1.039539594252776),
    @(14, 2, 1.02),
    @(14, 3, 1.031594923719023),
    @(14, 4, 1.034605596736237),
    @(14, 5, 1.039742949579308),
    @(14, 6, 1.047148197374006),
    @(14, 9, 1.033501694403902),
    @(14, 10, 1.038194910166922),
    @(14, 11, 1.038181603857067),
    @(14, 12, 1.043299938094984),
    @(14, 13, 1.050678134285065),
    @(14, 14, 1.03966926602198),
    @(15, 2, 1.02),
    @(15, 3, 1.031712545505114),
    @(15, 4, 1.034713753884612),
    @(15, 5, 1.039850890139605),
    @(15, 6, 1.04727413774109),
    @(15, 9, 1.03352323494899),
    @(15, 10, 1.038274666766518),
    @(15, 11, 1.038269831555892),
    @(15, 12, 1.04338806496281),
    @(15, 13, 1.050784359422255),
    @(15, 14, 1.039749135885095),
    @(16, 2, 1.02),
    @(16, 3, 1.032396966602232),
    @(16, 4, 1.035343198488974),
    @(16, 5, 1.040479152656704),
    @(16, 6, 1.048007267040619),
    @(16, 9, 1.033647746458606),
    @(16, 10, 1.038738500848666),
    @(16, 11, 1.038783053520977),
    @(16, 12, 1.043900777455332),
    @(16, 13, 1.051402526192915),
    @(16, 14, 1.040213628664843),
    @(17, 2, 1.02),
    @(17, 3, 1.032826120798568),
    @(17, 4, 1.035737965275437),
    @(17, 5, 1.040873247138212),
    @(17, 6, 1.048467230189275),
    @(17, 9, 1.033725088368804),
    @(17, 10, 1.039029114048861),
    @(17, 11, 1.039104719043395),
    @(17, 12, 1.044222191606372),
    @(17, 13, 1.051790190070824),
    @(17, 14, 1.040504654569115),
    @(18, 2, 1.02),
    @(18, 3, 1.033076380071585),
    @(18, 4, 1.035968202009085),
    @(18, 5, 1.041103116663925),
    @(18, 6, 1.048735551815086),
    @(18, 9, 1.033769926348942),
    @(18, 10, 1.039198502264008),
    @(18, 11, 1.039292245724476),
    @(18, 12, 1.044409596259903),
    @(18, 13, 1.052016273552025),
    @(18, 14, 1.040674283334957),
    @(19, 2, 1.02),
    @(19, 3, 1.033161702026527),
    @(19, 4, 1.036046702720289),
    @(19, 5, 1.041181496396535),
    @(19, 6, 1.048827048163282),
    @(19, 9, 1.033785168420345),
    @(19, 10, 1.039256238664711),
    @(19, 11, 1.039356171312957),
    @(19, 12, 1.044473484412731),
    @(19, 13, 1.052093356380089),
    @(19, 14, 1.040732101727971),
    @(20, 2, 1.02),
    @(20, 3, 1.032780082737035),
    @(20, 4, 1.035695613026663),
    @(20, 5, 1.040830964442186),
    @(20, 6, 1.048417877099333),
    @(20, 9, 1.033716818679502),
    @(20, 10, 1.038997946578088),
    @(20, 11, 1.039070217245359),
    @(20, 12, 1.044187714234603),
    @(20, 13, 1.05174860095825),
    @(20, 14, 1.040473442836959),
    @(21, 2, 1.02),
    @(21, 3, 1.031538704130509),
    @(21, 4, 1.034553902665801),
    @(21, 5, 1.039691360412084),
    @(21, 6, 1.047088007129042),
    @(21, 9, 1.033491384048459),
    @(21, 10, 1.038156784454996),
    @(21, 11, 1.038139430916808),
    @(21, 12, 1.043257814712395),
    @(21, 13, 1.050627363021781),
    @(21, 14, 1.03963108616717),
    @(22, 2, 1.02),
    @(22, 3, 1.030757659407213),
    @(22, 4, 1.033835843923209),
    @(22, 5, 1.038974850611133),
    @(22, 6, 1.046252155379431),
    @(22, 9, 1.033347174505263),
    @(22, 10, 1.037626814350042),
    @(22, 11, 1.037553346192891),
    @(22, 12, 1.042672508860471),
    @(22, 13, 1.049922082331765),
    @(22, 14, 1.039100363443878),
    @(23, 2, 1.02),
    @(23, 3, 1.031171757647999),
    @(23, 4, 1.034216521482721),
    @(23, 5, 1.039354685758833),
    @(23, 6, 1.046695228991362),
    @(23, 9, 1.033423856786612),
    @(23, 10, 1.037907865768087),
    @(23, 11, 1.037864122697908),
    @(23, 12, 1.042982851506013),
    @(23, 13, 1.050295994963216),
    @(23, 14, 1.039381813987174),
    @(24, 2, 1.02),
    @(24, 3, 1.032800885537581),
    @(24, 4, 1.035714750258905),
    @(24, 5, 1.040850070169494),
    @(24, 6, 1.048440177534284),
    @(24, 9, 1.033720556243429),
    @(24, 10, 1.03901203019186),
    @(24, 11, 1.039085807417222),
    @(24, 12, 1.044203293293301),
    @(24, 13, 1.051767393393714),
    @(24, 14, 1.040487546451078),
    @(25, 2, 1.02),
    @(25, 3, 1.03468908306846),
    @(25, 4, 1.037452419309767),
    @(25, 5, 1.042585416729144),
    @(25, 6, 1.050466386879072),
    @(25, 9, 1.034054100648774),
    @(25, 10, 1.040288584849325),
    @(25, 11, 1.040499764548423),
    @(25, 12, 1.045616772962991),
    @(25, 13, 1.053473526170743),
    @(25, 14, 1.04176591396258)
)

foreach ($entry in $updates) {
    $ws.Cells.Item($entry[0], $entry[1]).Value = $entry[2]
}
